# Calibrate VA transp, updated cap retirements
#
# Updates the "SoCDTtiNTY-psgr" sheet's calibration values (LDVs and
# motorbikes rows) and switches the active/selected sheet & cell to match
# the author's final view state (tab moves from "About" to
# "SoCDTtiNTY-psgr", selection lands on E15).

$wb = $excel.ActiveWorkbook

$psgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# LDVs row (row 2): calibrated down to a flat 7.5% across all years.
$psgr.Range("B2:H2").Value = 0.075

# motorbikes row (row 7): calibrated up to a flat 5.87% across all years.
$psgr.Range("B7:H7").Value = 0.0587

# Make the passenger sheet the active tab/selection, matching the saved
# view state in the authored workbook.
[void]$psgr.Activate()
[void]$psgr.Range("E15").Select()
